{"js": "const replacements = [\n  [\"2025-03-19 Wednesday\", \"2025-03-20 Thursday\"],\n  [\"477\u00d79=4293\", \"909\u00d74=3636\"],\n  [\"472\u00d77=3304\", \"626\u00d75=3130\"],\n  [\"591\u00d77=4137\", \"690\u00d72=1380\"],\n  [\"987\u00d78=7896\", \"435\u00d72=870\"],\n  [\"652\u00d74=2608\", \"392\u00d77=2744\"],\n  [\"511\u00d77=3577\", \"482\u00d72=964\"],\n  [\"354\u00d77=2478\", \"825\u00d77=5775\"],\n  [\"252\u00d77=1764\", \"151\u00d73=453\"],\n  [\"959\u00d73=2877\", \"832\u00d76=4992\"],\n  [\"148\u00d78=1184\", \"563\u00d76=3378\"],\n  [\"641\u00d78=5128\", \"726\u00d74=2904\"],\n  [\"232\u00d78=1856\", \"470\u00d77=3290\"],\n  [\"322\u00d74=1288\", \"674\u00d73=2022\"],\n  [\"689\u00d72=1378\", \"441\u00d72=882\"],\n  [\"498\u00d79=4482\", \"341\u00d79=3069\"],\n  [\"889\u00d75=4445\", \"931\u00d79=8379\"],\n  [\"212\u00d76=1272\", \"358\u00d78=2864\"],\n  [\"626\u00d77=4382\", \"797\u00d79=7173\"],\n  [\"222\u00d74=888\", \"743\u00d79=6687\"],\n  [\"114\u00d77=798\", \"217\u00d76=1302\"],\n  [\"820\u00d73=2460\", \"427\u00d75=2135\"],\n  [\"596\u00d76=3576\", \"420\u00d77=2940\"],\n  [\"694\u00d76=4164\", \"247\u00d78=1976\"],\n  [\"284\u00d78=2272\", \"269\u00d73=807\"],\n  [\"637\u00d72=1274\", \"810\u00d75=4050\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  // Each snippet is unique in this document, so replace every match found\n  // (normally exactly one) while keeping the run's existing formatting.\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\nreturn \"done\";\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-03-19 Wednesday\", \"2025-03-20 Thursday\"),\n    @(\"477\u00d79=4293\", \"909\u00d74=3636\"),\n    @(\"472\u00d77=3304\", \"626\u00d75=3130\"),\n    @(\"591\u00d77=4137\", \"690\u00d72=1380\"),\n    @(\"987\u00d78=7896\", \"435\u00d72=870\"),\n    @(\"652\u00d74=2608\", \"392\u00d77=2744\"),\n    @(\"511\u00d77=3577\", \"482\u00d72=964\"),\n    @(\"354\u00d77=2478\", \"825\u00d77=5775\"),\n    @(\"252\u00d77=1764\", \"151\u00d73=453\"),\n    @(\"959\u00d73=2877\", \"832\u00d76=4992\"),\n    @(\"148\u00d78=1184\", \"563\u00d76=3378\"),\n    @(\"641\u00d78=5128\", \"726\u00d74=2904\"),\n    @(\"232\u00d78=1856\", \"470\u00d77=3290\"),\n    @(\"322\u00d74=1288\", \"674\u00d73=2022\"),\n    @(\"689\u00d72=1378\", \"441\u00d72=882\"),\n    @(\"498\u00d79=4482\", \"341\u00d79=3069\"),\n    @(\"889\u00d75=4445\", \"931\u00d79=8379\"),\n    @(\"212\u00d76=1272\", \"358\u00d78=2864\"),\n    @(\"626\u00d77=4382\", \"797\u00d79=7173\"),\n    @(\"222\u00d74=888\", \"743\u00d79=6687\"),\n    @(\"114\u00d77=798\", \"217\u00d76=1302\"),\n    @(\"820\u00d73=2460\", \"427\u00d75=2135\"),\n    @(\"596\u00d76=3576\", \"420\u00d77=2940\"),\n    @(\"694\u00d76=4164\", \"247\u00d78=1976\"),\n    @(\"284\u00d78=2272\", \"269\u00d73=807\"),\n    @(\"637\u00d72=1274\", \"810\u00d75=4050\"),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n\nWrite-Output \"done\"\n"}
